$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the quantity for the "2x3 header" BOM line (row 3): 2 -> 1.
# G3 (Total Cost = F3*E3) recalculates automatically.
$ws.Range("F3").Value = 1

# Restore the last-known cursor position/selection to F4.
$ws.Range("F4").Select() | Out-Null
